$wb = $excel.ActiveWorkbook

# --- "Inspire" sheet (sheet4.xml) ---
$inspire = $wb.Worksheets.Item("Inspire")

# New row 13
$inspire.Range("B13").Value = "random distort"

# New row 14
$inspire.Range("B14").Value = "anti training"
$inspire.Range("C14").Value = "existing models to transfer learning"

$inspire.Range("C13").Value = "to many similar inputs,  and average classifier"

# Row 12: update the existing "!!!!!" note to add more detail
$inspire.Range("C12").Value = "!!!!!  Use this example to train."

$null = $inspire.Range("C12").Select()

# --- "Todo" sheet (sheet5.xml) ---
$todo = $wb.Worksheets.Item("Todo")

# New column C needs a wider width to hold the new notes
$todo.Columns.Item(3).ColumnWidth = 34

# Row 3: tweak existing wording
$todo.Range("B3").Value = "input resize keep ratio"

# New row 4
$todo.Range("B4").Value = "get adv training datasets"
$todo.Range("C4").Value = "augment,  MIM output, MSB"

$null = $todo.Range("E13").Select()
